$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"24.66000000000042"
$ws.Range("H2").Value = [double]"4.023298361133243e-10"
$ws.Range("I2").Value = [double]"4.023298361133243e-10"
$ws.Range("L2").Value = [double]"59.74673597607177"
$ws.Range("M2").Value = "[43.638906444910546, 75.85456550723299]"
$ws.Range("N2").Value = [double]"2.053646808164444e-09"
$ws.Range("O2").Value = [double]"2.053646808164444e-09"
$ws.Range("P2").Value = [double]"1.427710775505271"
$ws.Range("Q2").Value = "[1.1006580868432705, 1.7547634641672714]"
$ws.Range("R2").Value = [double]"2.495070816621592e-11"
$ws.Range("S2").Value = [double]"2.495070816621592e-11"
$ws.Range("T2").Value = [double]"58.18703974034965"
$ws.Range("U2").Value = "[48.103927465537, 68.2701520151623]"
$ws.Range("V2").Value = [double]"3.774758283725532e-15"
$ws.Range("W2").Value = [double]"3.774758283725532e-15"
$ws.Range("X2").Value = [double]"19.0565765765769"
$ws.Range("Y2").Value = [double]"17.77297297297328"
$ws.Range("Z2").Value = [double]"20.34018018018052"
$ws.Range("F3").Value = [double]"24.66000000000042"
$ws.Range("H3").Value = [double]"3.026022765695302e-10"
$ws.Range("I3").Value = [double]"3.026022765695302e-10"
$ws.Range("L3").Value = [double]"55.42544686098345"
$ws.Range("M3").Value = "[38.868251650335864, 71.98264207163103]"
$ws.Range("N3").Value = [double]"2.459567549628616e-08"
$ws.Range("O3").Value = [double]"2.459567549628616e-08"
$ws.Range("P3").Value = [double]"1.276763380738194"
$ws.Range("Q3").Value = "[0.9497106920761933, 1.6038160694001942]"
$ws.Range("R3").Value = [double]"5.464171337621337e-10"
$ws.Range("S3").Value = [double]"5.464171337621337e-10"
$ws.Range("T3").Value = [double]"52.63079801481621"
$ws.Range("U3").Value = "[43.38540192342182, 61.8761941062106]"
$ws.Range("V3").Value = [double]"5.995204332975845e-15"
$ws.Range("W3").Value = [double]"5.995204332975845e-15"
$ws.Range("X3").Value = [double]"19.64900900900934"
$ws.Range("Y3").Value = [double]"18.36540540540572"
$ws.Range("Z3").Value = [double]"20.93261261261296"
$ws.Range("F4").Value = [double]"24.66000000000042"
$ws.Range("H4").Value = [double]"2.353276462585541e-10"
$ws.Range("I4").Value = [double]"2.353276462585541e-10"
$ws.Range("L4").Value = [double]"65.69547765348938"
$ws.Range("M4").Value = "[44.71329352292513, 86.67766178405363]"
$ws.Range("N4").Value = [double]"1.095031156861381e-07"
$ws.Range("O4").Value = [double]"1.095031156861381e-07"
$ws.Range("P4").Value = [double]"0.8993948938205003"
$ws.Range("Q4").Value = "[0.5849211547224238, 1.2138686329185768]"
$ws.Range("R4").Value = [double]"7.074636645398868e-07"
$ws.Range("S4").Value = [double]"7.074636645398868e-07"
$ws.Range("T4").Value = [double]"65.98063980146298"
$ws.Range("U4").Value = "[55.220534707790776, 76.7407448951352]"
$ws.Range("V4").Value = [double]"4.440892098500626e-16"
$ws.Range("W4").Value = [double]"4.440892098500626e-16"
$ws.Range("X4").Value = [double]"21.13009009009045"
$ws.Range("Y4").Value = [double]"19.8958558558562"
$ws.Range("Z4").Value = [double]"22.3643243243247"
$ws.Range("F5").Value = [double]"24.66000000000042"
$ws.Range("H5").Value = [double]"1.085830814151478e-07"
$ws.Range("I5").Value = [double]"1.085830814151478e-07"
$ws.Range("L5").Value = [double]"51.13213767676623"
$ws.Range("M5").Value = "[30.56854372665842, 71.69573162687405]"
$ws.Range("N5").Value = [double]"8.952991988886438e-06"
$ws.Range("O5").Value = [double]"8.952991988886438e-06"
$ws.Range("P5").Value = [double]"1.062921238151501"
$ws.Range("Q5").Value = "[0.6603948521059619, 1.4654476241970391]"
$ws.Range("R5").Value = [double]"3.16401645594766e-06"
$ws.Range("S5").Value = [double]"3.16401645594766e-06"
$ws.Range("T5").Value = [double]"58.18915496383845"
$ws.Range("U5").Value = "[47.51156226079952, 68.86674766687739]"
$ws.Range("V5").Value = [double]"2.597921877622866e-14"
$ws.Range("W5").Value = [double]"2.597921877622866e-14"
$ws.Range("X5").Value = [double]"20.48828828828864"
$ws.Range("Y5").Value = [double]"18.90846846846879"
$ws.Range("Z5").Value = [double]"22.06810810810848"
$ws.Range("F6").Value = [double]"25.49000000000055"
$ws.Range("H6").Value = [double]"4.199117109582318e-08"
$ws.Range("I6").Value = [double]"4.199117109582318e-08"
$ws.Range("L6").Value = [double]"57.53962757335051"
$ws.Range("M6").Value = "[36.238744523237074, 78.84051062346396]"
$ws.Range("N6").Value = [double]"2.094620282555582e-06"
$ws.Range("O6").Value = [double]"2.094620282555582e-06"
$ws.Range("P6").Value = [double]"0.5597632555945768"
$ws.Range("Q6").Value = "[0.18239476867688342, 0.9371317425122703]"
$ws.Range("R6").Value = [double]"0.004540609241171945"
$ws.Range("S6").Value = [double]"0.004540609241171945"
$ws.Range("T6").Value = [double]"60.44361050961876"
$ws.Range("U6").Value = "[49.22562931577117, 71.66159170346634]"
$ws.Range("V6").Value = [double]"3.774758283725532e-14"
$ws.Range("W6").Value = [double]"3.774758283725532e-14"
$ws.Range("X6").Value = [double]"23.21911911911962"
$ws.Range("Y6").Value = [double]"21.68818818818865"
$ws.Range("Z6").Value = [double]"24.75005005005059"
$ws.Range("F7").Value = [double]"25.49000000000055"
$ws.Range("H7").Value = [double]"7.400802193302525e-11"
$ws.Range("I7").Value = [double]"7.400802193302525e-11"
$ws.Range("L7").Value = [double]"64.98475038205004"
$ws.Range("M7").Value = "[48.85093237757613, 81.11856838652395]"
$ws.Range("N7").Value = [double]"2.36555441901487e-10"
$ws.Range("O7").Value = [double]"2.36555441901487e-10"
$ws.Range("P7").Value = [double]"0.1069210712933462"
$ws.Range("Q7").Value = "[-0.16981581911296217, 0.3836579616996545]"
$ws.Range("R7").Value = [double]"0.4405376480218617"
$ws.Range("S7").Value = [double]"0.4405376480218617"
$ws.Range("T7").Value = [double]"55.35028418626226"
$ws.Range("U7").Value = "[45.43766647979457, 65.26290189272994]"
$ws.Range("V7").Value = [double]"1.154631945610163e-14"
$ws.Range("W7").Value = [double]"1.154631945610163e-14"
$ws.Range("X7").Value = [double]"25.05623623623677"
$ws.Range("Y7").Value = [double]"23.93355355355407"
$ws.Range("Z7").Value = [double]"26.17891891891948"
$ws.Range("F8").Value = [double]"25.49000000000055"
$ws.Range("H8").Value = [double]"1.180819220181206e-09"
$ws.Range("I8").Value = [double]"1.180819220181206e-09"
$ws.Range("L8").Value = [double]"59.78535539239688"
$ws.Range("M8").Value = "[44.21675214186699, 75.35395864292677]"
$ws.Range("N8").Value = [double]"8.41796854444965e-10"
$ws.Range("O8").Value = [double]"8.41796854444965e-10"
$ws.Range("P8").Value = [double]"-0.006289474781961957"
$ws.Range("Q8").Value = "[-0.30818426431611634, 0.2956053147521924]"
$ws.Range("R8").Value = [double]"0.9667158264814162"
$ws.Range("S8").Value = [double]"0.9667158264814162"
$ws.Range("T8").Value = [double]"56.33230225590744"
$ws.Range("U8").Value = "[46.238253363079856, 66.42635114873502]"
$ws.Range("V8").Value = [double]"1.176836406102666e-14"
$ws.Range("W8").Value = [double]"1.176836406102666e-14"
$ws.Range("X8").Value = [double]"0.02551551551551512"
$ws.Range("Y8").Value = [double]"-1.199229229229256"
$ws.Range("Z8").Value = [double]"1.250260260260287"
$ws.Range("F9").Value = [double]"25.49000000000055"
$ws.Range("H9").Value = [double]"1.157869744528028e-08"
$ws.Range("I9").Value = [double]"1.157869744528028e-08"
$ws.Range("L9").Value = [double]"59.73762502072056"
$ws.Range("M9").Value = "[37.763333475899685, 81.71191656554144]"
$ws.Range("N9").Value = [double]"1.86235323140771e-06"
$ws.Range("O9").Value = [double]"1.86235323140771e-06"
$ws.Range("P9").Value = [double]"-0.6289474781961548"
$ws.Range("Q9").Value = "[-0.9937370155499243, -0.2641579408423853]"
$ws.Range("R9").Value = [double]"0.001150453559177755"
$ws.Range("S9").Value = [double]"0.001150453559177755"
$ws.Range("T9").Value = [double]"66.5788295605233"
$ws.Range("U9").Value = "[55.17130446195645, 77.98635465909015]"
$ws.Range("V9").Value = [double]"2.664535259100376e-15"
$ws.Range("W9").Value = [double]"2.664535259100376e-15"
$ws.Range("X9").Value = [double]"2.551551551551608"
$ws.Range("Y9").Value = [double]"1.071651651651677"
$ws.Range("Z9").Value = [double]"4.031451451451538"
$ws.Range("F10").Value = [double]"25.49000000000055"
$ws.Range("H10").Value = [double]"1.03755892766344e-10"
$ws.Range("I10").Value = [double]"1.03755892766344e-10"
$ws.Range("L10").Value = [double]"59.77798765287794"
$ws.Range("M10").Value = "[41.10968081849647, 78.44629448725941]"
$ws.Range("N10").Value = [double]"6.704940114232727e-08"
$ws.Range("O10").Value = [double]"6.704940114232727e-08"
$ws.Range("P10").Value = [double]"-0.9182633181663862"
$ws.Range("Q10").Value = "[-1.2201581077005406, -0.6163685286322318]"
$ws.Range("R10").Value = [double]"2.027648813296423e-07"
$ws.Range("S10").Value = [double]"2.027648813296423e-07"
$ws.Range("T10").Value = [double]"53.95033514866888"
$ws.Range("U10").Value = "[44.20544254796524, 63.69522774937251]"
$ws.Range("V10").Value = [double]"1.532107773982716e-14"
$ws.Range("W10").Value = [double]"1.532107773982716e-14"
$ws.Range("X10").Value = [double]"3.725265265265346"
$ws.Range("Y10").Value = [double]"2.500520520520574"
$ws.Range("Z10").Value = [double]"4.950010010010118"
